$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.008.35"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "3.411.20"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.41"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.65"
$ws.Range("E6").Value = "  -4.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("E7").Value = "  +7.10%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.730"
$ws.Range("E9").Value = "  +7.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  +17.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.38"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000219"
$ws.Range("E12").Value = "  +67.77%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "3.952.93"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.92"
$ws.Range("E15").Value = "  +5.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.80"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "3.396.47"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.05"
$ws.Range("E18").Value = "  +9.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.07"
$ws.Range("E19").Value = "  +5.04%  "
$ws.Range("D20").Value = "61.932.31"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "406.64"
$ws.Range("E21").Value = "  +29.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.52"
$ws.Range("E22").Value = "  +5.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.18"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.08"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.24"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "32.88"
$ws.Range("E26").Value = "  +10.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.67"
$ws.Range("E27").Value = "  +5.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.86"
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.93"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0498"
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.11"
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("E40").Value = "  +6.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.311"
$ws.Range("E42").Value = "  +4.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.77"
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.08"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("E46").Value = "  +8.59%  "
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.95"
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("D49").Value = "2.117.25"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  +4.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.131"
$ws.Range("E51").Value = "  +15.17%  "
